$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update F-column query timestamps on the data sheet ("time_taken")
$ws1.Cells.Item(2,6).Value = "2021-10-05 14:21:36.351780"
$ws1.Cells.Item(3,6).Value = "2021-10-05 14:21:36.351788"
$ws1.Cells.Item(4,6).Value = "2021-10-05 14:21:36.351791"
$ws1.Cells.Item(5,6).Value = "2021-10-05 14:21:36.351794"
$ws1.Cells.Item(6,6).Value = "2021-10-05 14:21:36.351797"
$ws1.Cells.Item(7,6).Value = "2021-10-05 14:21:36.351799"
$ws1.Cells.Item(8,6).Value = "2021-10-05 14:21:36.351802"
$ws1.Cells.Item(9,6).Value = "2021-10-05 14:21:36.351804"
$ws1.Cells.Item(10,6).Value = "2021-10-05 14:21:36.351807"
$ws1.Cells.Item(11,6).Value = "2021-10-05 14:21:36.351810"
$ws1.Cells.Item(12,6).Value = "2021-10-05 14:21:36.351812"
$ws1.Cells.Item(13,6).Value = "2021-10-05 14:21:36.351815"
$ws1.Cells.Item(14,6).Value = "2021-10-05 14:21:36.351817"
$ws1.Cells.Item(15,6).Value = "2021-10-05 14:21:36.351820"
$ws1.Cells.Item(16,6).Value = "2021-10-05 14:21:36.351822"
$ws1.Cells.Item(17,6).Value = "2021-10-05 14:21:36.351824"
$ws1.Cells.Item(18,6).Value = "2021-10-05 14:21:36.351827"
$ws1.Cells.Item(19,6).Value = "2021-10-05 14:21:36.351830"
$ws1.Cells.Item(20,6).Value = "2021-10-05 14:21:36.351832"
$ws1.Cells.Item(21,6).Value = "2021-10-05 14:21:36.351835"
$ws1.Cells.Item(22,6).Value = "2021-10-05 14:21:36.351837"
$ws1.Cells.Item(23,6).Value = "2021-10-05 14:21:36.351840"
$ws1.Cells.Item(24,6).Value = "2021-10-05 14:21:36.351842"
$ws1.Cells.Item(25,6).Value = "2021-10-05 14:21:36.351845"
$ws1.Cells.Item(26,6).Value = "2021-10-05 14:21:36.351847"
$ws1.Cells.Item(27,6).Value = "2021-10-05 14:21:36.351850"
$ws1.Cells.Item(28,6).Value = "2021-10-05 14:21:36.351852"
$ws1.Cells.Item(29,6).Value = "2021-10-05 14:21:36.351855"
$ws1.Cells.Item(30,6).Value = "2021-10-05 14:21:36.351857"
$ws1.Cells.Item(31,6).Value = "2021-10-05 14:21:36.351860"
$ws1.Cells.Item(32,6).Value = "2021-10-05 14:21:36.351862"
$ws1.Cells.Item(33,6).Value = "2021-10-05 14:21:36.351865"
$ws1.Cells.Item(34,6).Value = "2021-10-05 14:21:36.351867"
$ws1.Cells.Item(35,6).Value = "2021-10-05 14:21:36.351870"
$ws1.Cells.Item(36,6).Value = "2021-10-05 14:21:36.351873"
$ws1.Cells.Item(37,6).Value = "2021-10-05 14:21:36.351875"
$ws1.Cells.Item(38,6).Value = "2021-10-05 14:21:36.351877"
$ws1.Cells.Item(39,6).Value = "2021-10-05 14:21:36.351880"
$ws1.Cells.Item(40,6).Value = "2021-10-05 14:21:36.351882"
$ws1.Cells.Item(41,6).Value = "2021-10-05 14:21:36.351885"
$ws1.Cells.Item(42,6).Value = "2021-10-05 14:21:36.351888"
$ws1.Cells.Item(43,6).Value = "2021-10-05 14:21:36.351890"
$ws1.Cells.Item(44,6).Value = "2021-10-05 14:21:36.351893"
$ws1.Cells.Item(45,6).Value = "2021-10-05 14:21:36.351895"
$ws1.Cells.Item(46,6).Value = "2021-10-05 14:21:36.351898"
$ws1.Cells.Item(47,6).Value = "2021-10-05 14:21:36.351900"
$ws1.Cells.Item(48,6).Value = "2021-10-05 14:21:36.351903"
$ws1.Cells.Item(49,6).Value = "2021-10-05 14:21:36.351905"
$ws1.Cells.Item(50,6).Value = "2021-10-05 14:21:36.351908"
$ws1.Cells.Item(51,6).Value = "2021-10-05 14:21:36.351910"
$ws1.Cells.Item(52,6).Value = "2021-10-05 14:21:36.351912"
$ws1.Cells.Item(53,6).Value = "2021-10-05 14:21:36.351915"
$ws1.Cells.Item(54,6).Value = "2021-10-05 14:21:36.351918"
$ws1.Cells.Item(55,6).Value = "2021-10-05 14:21:36.351920"
$ws1.Cells.Item(56,6).Value = "2021-10-05 14:21:36.351923"
$ws1.Cells.Item(57,6).Value = "2021-10-05 14:21:36.351925"
$ws1.Cells.Item(58,6).Value = "2021-10-05 14:21:36.351928"
$ws1.Cells.Item(59,6).Value = "2021-10-05 14:21:36.351930"
$ws1.Cells.Item(60,6).Value = "2021-10-05 14:21:36.351933"
$ws1.Cells.Item(61,6).Value = "2021-10-05 14:21:36.351935"
$ws1.Cells.Item(62,6).Value = "2021-10-05 14:21:36.351938"
$ws1.Cells.Item(63,6).Value = "2021-10-05 14:21:36.351940"
$ws1.Cells.Item(64,6).Value = "2021-10-05 14:21:36.351943"
$ws1.Cells.Item(65,6).Value = "2021-10-05 14:21:36.351945"
$ws1.Cells.Item(66,6).Value = "2021-10-05 14:21:36.351949"
$ws1.Cells.Item(67,6).Value = "2021-10-05 14:21:36.351951"
$ws1.Cells.Item(68,6).Value = "2021-10-05 14:21:36.351954"
$ws1.Cells.Item(69,6).Value = "2021-10-05 14:21:36.351956"
$ws1.Cells.Item(70,6).Value = "2021-10-05 14:21:36.351959"
$ws1.Cells.Item(71,6).Value = "2021-10-05 14:21:36.351961"
$ws1.Cells.Item(72,6).Value = "2021-10-05 14:21:36.351964"
$ws1.Cells.Item(73,6).Value = "2021-10-05 14:21:36.351966"
$ws1.Cells.Item(74,6).Value = "2021-10-05 14:21:36.351969"
$ws1.Cells.Item(75,6).Value = "2021-10-05 14:21:36.351971"
$ws1.Cells.Item(76,6).Value = "2021-10-05 14:21:36.351974"
$ws1.Cells.Item(77,6).Value = "2021-10-05 14:21:36.351976"
$ws1.Cells.Item(78,6).Value = "2021-10-05 14:21:36.351981"

# Add the new "metadata" worksheet positioned after "data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Reuse the header-row style (bold + border + centered) from the data sheet
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

# Reuse the index-column style from the data sheet for A2
$ws1.Cells.Item(2,1).Copy()
$ws2.Cells.Item(2,1).PasteSpecial(-4122)

# Header row
$ws2.Cells.Item(1,2).Value = "data_name"
$ws2.Cells.Item(1,3).Value = "data_id"
$ws2.Cells.Item(1,4).Value = "data_version"
$ws2.Cells.Item(1,5).Value = "data_version_created"
$ws2.Cells.Item(1,6).Value = "panel_query_time"
$ws2.Cells.Item(1,7).Value = "panel_get_request"

# Data row
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "Monogenic diabetes"
$ws2.Cells.Item(2,3).Value = 472
$ws2.Cells.Item(2,4).NumberFormat = "@"
$ws2.Cells.Item(2,4).Value = "2.43"
$ws2.Cells.Item(2,5).Value = "2021-07-28T09:59:13.775656Z"
$ws2.Cells.Item(2,6).Value = "2021-10-05 14:21:36.348091"
$ws2.Cells.Item(2,7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/472/?format=json"

Write-Output "Applied metadata sheet and refreshed time_taken timestamps"
